# Generate Report for Archive
# The handoff has moved on: update the status text everywhere it is
# shown (Overview rollup + each per-locale status table), then let the
# status columns re-size to fit the new (shorter) text, same as Excel
# does automatically after a content refresh.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = "In Translation"

# --- Per-locale detail sheets: Status column (C2) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Re-fit the columns that held the old, longer "Ready for handoff"
# text now that the shorter "In Translation" text is in place.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
